# CommonTypos.docx edit:
#  1. Insert 8 new bullet/list paragraphs ("Today I did an interview..." ...
#     "Focus on pushing out all the code...") right before the
#     "Interview Improvements:" heading paragraph.
#  2. Add a <w:lastRenderedPageBreak/> run-child at the start of the run
#     that holds "However, I couldn't get the in-place solution...".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: insert the eight new list paragraphs before "Interview
# Improvements:"
# ---------------------------------------------------------------------

$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -eq "Interview Improvements:`r") {
        $targetIdx = $idx
    }
}

if ($targetIdx -eq -1) {
    throw "Could not find 'Interview Improvements:' paragraph"
}

$target = $d.Paragraphs($targetIdx)
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)

# NOTE: a pkg:package / pkg:part wrapper is required for InsertXML to
# splice content in (rather than stomp the surrounding paragraph), and
# a trailing empty <w:p/> is required because the *last* <w:p> in an
# InsertXML payload never becomes its own paragraph -- its content (here,
# none) is merged onto the paragraph at the insertion point instead of
# breaking before it. Without the dummy trailing <w:p/> our real eighth
# paragraph ("Focus on pushing out...") would lose its own paragraph
# mark/formatting and get merged into "Interview Improvements:".
$insertXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Today I did an interview that was literally just program a game in 3 hours</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>I was TOTALLY unprepared for this</w:t></w:r><w:r><w:t>, but did my best to wing it in 3 hours</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Make sure to ask the interviewer about the exact format for the interview</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Frankly, if I was better at Python 3 programming I could have actually done this</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>However, I am straight up not good enough at coding to work that quickly under pressure</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I'm not going to specifically prepare for interviews like this, the only way to get better is simply to code more </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The actual plan here is to implement a Django Web application over Winter Break (Mid-December to February 1</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Focus on pushing out all the code and learning as much about Python as possible</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
# swap in the proper curly apostrophe (kept out of the single-quoted
# here-string above to dodge quoting issues)
$insertXml = $insertXml.Replace("I'm not going", "I" + [char]0x2019 + "m not going")

$insertionPoint.InsertXML($insertXml) | Out-Null

# Re-locate & drop the dummy empty paragraph the trailing <w:p/> left
# behind right before "Interview Improvements:".
$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -eq "Interview Improvements:`r") {
        $targetIdx = $idx
    }
}
$dummy = $d.Paragraphs($targetIdx - 1)
if ($dummy.Range.Text -eq "`r") {
    $dummy.Range.Delete()
}

# ---------------------------------------------------------------------
# Part 2: add <w:lastRenderedPageBreak/> to the run with "However, I
# couldn't get the in-place solution..."
# ---------------------------------------------------------------------

$breakIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*in-place solution*") {
        $breakIdx = $idx
    }
}

if ($breakIdx -eq -1) {
    throw "Could not find 'in-place solution' paragraph"
}

$bp = $d.Paragraphs($breakIdx)
$bpText = $bp.Range.Text
$bpText = $bpText.Substring(0, $bpText.Length - 1)  # drop trailing paragraph mark
$bpRange = $d.Range($bp.Range.Start, $bp.Range.End - 1)

$runXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>' + $bpText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$bpRange.InsertXML($runXml) | Out-Null

Write-Output "edit complete"
